$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I1").Value = 'Other found locations'

$ws.Range("E2").Value = '[Harrison X.%Bai%NULL%2,       Ben%Hsieh%NULL%1,       Zeng%Xiong%NULL%1,       Kasey%Halsey%NULL%1,       Ji Whae%Choi%NULL%1,       Thi My Linh%Tran%NULL%1,       Ian%Pan%NULL%1,       Lin-Bo%Shi%NULL%1,       Dong-Cui%Wang%NULL%1,       Ji%Mei%NULL%1,       Xiao-Long%Jiang%NULL%1,       Qiu-Hua%Zeng%NULL%1,       Thomas K.%Egglin%NULL%1,       Ping-Feng%Hu%NULL%1,       Saurabh%Agarwal%NULL%1,       Fangfang%Xie%NULL%1,       Sha%Li%NULL%1,       Terrance%Healey%NULL%1,       Michael K.%Atalay%NULL%1,       Wei-Hua%Liao%liaoweihua2017@163.com%1,      Harrison X.%Bai%NULL%1,      Ben%Hsieh%NULL%1,      Zeng%Xiong%NULL%1,      Kasey%Halsey%NULL%1,      Ji Whae%Choi%NULL%1,      Thi My Linh%Tran%NULL%1,      Ian%Pan%NULL%1,      Lin-Bo%Shi%NULL%1,      Dong-Cui%Wang%NULL%1,      Ji%Mei%NULL%1,      Xiao-Long%Jiang%NULL%1,      Qiu-Hua%Zeng%NULL%1,      Thomas K.%Egglin%NULL%1,      Ping-Feng%Hu%NULL%1,      Saurabh%Agarwal%NULL%1,      Fangfang%Xie%NULL%1,      Sha%Li%NULL%1,      Terrance%Healey%NULL%1,      Michael K.%Atalay%NULL%1,      Wei-Hua%Liao%liaoweihua2017@163.com%1]'
$ws.Range("I2").Value = '_PMC'

$ws.Range("I3").Value = ''

$ws.Range("E4").Value = '[Damiano%Caruso%NULL%2,       Marta%Zerunian%NULL%1,       Michela%Polici%NULL%1,       Francesco%Pucciarelli%NULL%1,       Tiziano%Polidori%NULL%1,       Carlotta%Rucci%NULL%1,       Gisella%Guido%NULL%1,       Benedetta%Bracci%NULL%1,       Chiara%de Dominicis%NULL%1,       Andrea%Laghi%andrea.laghi@uniroma1.it%1,      Damiano%Caruso%NULL%1,      Marta%Zerunian%NULL%1,      Michela%Polici%NULL%1,      Francesco%Pucciarelli%NULL%1,      Tiziano%Polidori%NULL%1,      Carlotta%Rucci%NULL%1,      Gisella%Guido%NULL%1,      Benedetta%Bracci%NULL%1,      Chiara%de Dominicis%NULL%1,      Andrea%Laghi%andrea.laghi@uniroma1.it%1]'
$ws.Range("I4").Value = '_PMC'

$ws.Range("E5").Value = '[Xiaofeng%Chen%NULL%2,       Yanyan%Tang%NULL%1,       Yongkang%Mo%NULL%1,       Shengkai%Li%NULL%1,       Daiying%Lin%NULL%1,       Zhijian%Yang%NULL%1,       Zhiqi%Yang%NULL%1,       Hongfu%Sun%NULL%1,       Jinming%Qiu%NULL%1,       Yuting%Liao%NULL%1,       Jianning%Xiao%NULL%1,       Xiangguang%Chen%NULL%1,       Xianheng%Wu%NULL%1,       Renhua%Wu%NULL%1,       Zhuozhi%Dai%zhuozhi@ualberta.ca%1,      Xiaofeng%Chen%NULL%1,      Yanyan%Tang%NULL%1,      Yongkang%Mo%NULL%1,      Shengkai%Li%NULL%1,      Daiying%Lin%NULL%1,      Zhijian%Yang%NULL%1,      Zhiqi%Yang%NULL%1,      Hongfu%Sun%NULL%1,      Jinming%Qiu%NULL%1,      Yuting%Liao%NULL%1,      Jianning%Xiao%NULL%1,      Xiangguang%Chen%NULL%1,      Xianheng%Wu%NULL%1,      Renhua%Wu%NULL%1,      Zhuozhi%Dai%zhuozhi@ualberta.ca%1]'
$ws.Range("I5").Value = '_PMC_Springer'

$ws.Range("I6").Value = ''

$ws.Range("E7").Value = '[Hyewon%Choi%NULL%2,       Xiaolong%Qi%NULL%1,       Soon Ho%Yoon%yshoka@gmail.com%1,       Sang Joon%Park%NULL%1,       Kyung Hee%Lee%NULL%1,       Jin Yong%Kim%NULL%1,       Young Kyung%Lee%NULL%1,       Hongseok%Ko%NULL%1,       Ki Hwan%Kim%NULL%1,       Chang Min%Park%NULL%1,       Yun-Hyeon%Kim%NULL%1,       Junqiang%Lei%NULL%1,       Jung Hee%Hong%NULL%1,       Hyungjin%Kim%NULL%1,       Eui Jin%Hwang%NULL%1,       Seung Jin%Yoo%NULL%1,       Ju Gang%Nam%NULL%1,       Chang Hyun%Lee%NULL%1,       Jin Mo%Goo%NULL%1,      Hyewon%Choi%NULL%1,      Xiaolong%Qi%NULL%1,      Soon Ho%Yoon%yshoka@gmail.com%1,      Sang Joon%Park%NULL%1,      Kyung Hee%Lee%NULL%1,      Jin Yong%Kim%NULL%1,      Young Kyung%Lee%NULL%1,      Hongseok%Ko%NULL%1,      Ki Hwan%Kim%NULL%1,      Chang Min%Park%NULL%1,      Yun-Hyeon%Kim%NULL%1,      Junqiang%Lei%NULL%1,      Jung Hee%Hong%NULL%1,      Hyungjin%Kim%NULL%1,      Eui Jin%Hwang%NULL%1,      Seung Jin%Yoo%NULL%1,      Ju Gang%Nam%NULL%1,      Chang Hyun%Lee%NULL%1,      Jin Mo%Goo%NULL%1]'
$ws.Range("I7").Value = '_PMC'

$ws.Range("E8").Value = '[Yuki%Himoto%yukihimoto@gmail.com%2,       Akihiko%Sakata%NULL%1,       Mitsuhiro%Kirita%NULL%1,       Takashi%Hiroi%NULL%1,       Ken-ichiro%Kobayashi%NULL%1,       Kenji%Kubo%NULL%1,       Hyunjin%Kim%NULL%1,       Azusa%Nishimoto%NULL%1,       Chikara%Maeda%NULL%1,       Akira%Kawamura%NULL%1,       Nobuhiro%Komiya%NULL%1,       Shigeaki%Umeoka%NULL%1,      Yuki%Himoto%yukihimoto@gmail.com%1,      Akihiko%Sakata%NULL%1,      Mitsuhiro%Kirita%NULL%1,      Takashi%Hiroi%NULL%1,      Ken-ichiro%Kobayashi%NULL%1,      Kenji%Kubo%NULL%1,      Hyunjin%Kim%NULL%1,      Azusa%Nishimoto%NULL%1,      Chikara%Maeda%NULL%1,      Akira%Kawamura%NULL%1,      Nobuhiro%Komiya%NULL%1,      Shigeaki%Umeoka%NULL%1]'
$ws.Range("I8").Value = '_PMC_Springer'

$ws.Range("E9").Value = '[Chunqin%Long%NULL%2,       Huaxiang%Xu%NULL%1,       Qinglin%Shen%NULL%1,       Xianghai%Zhang%NULL%1,       Bing%Fan%26171381@qq.com%1,       Chuanhong%Wang%NULL%1,       Bingliang%Zeng%NULL%1,       Zicong%Li%NULL%1,       Xiaofen%Li%NULL%1,       Honglu%Li%NULL%1,      Chunqin%Long%NULL%1,      Huaxiang%Xu%NULL%1,      Qinglin%Shen%NULL%1,      Xianghai%Zhang%NULL%1,      Bing%Fan%26171381@qq.com%1,      Chuanhong%Wang%NULL%1,      Bingliang%Zeng%NULL%1,      Zicong%Li%NULL%1,      Xiaofen%Li%NULL%1,      Honglu%Li%NULL%1]'
$ws.Range("I9").Value = '_PMC_elsevier'

$ws.Range("E10").Value = '[Anthony%Dangis%NULL%4,      Christopher%Gieraerts%NULL%1,      Yves%De Bruecker%NULL%1,      Lode%Janssen%NULL%1,      Hanne%Valgaeren%NULL%1,      Dagmar%Obbels%NULL%1,      Marc%Gillis%NULL%1,      Marc%Van Ranst%NULL%0,      Johan%Frans%NULL%1,      Annick%Demeyere%NULL%1,      Rolf%Symons%rolf.symons@imelda.be%1]'
$ws.Range("F10").Value = 'PMC7166306'
$ws.Range("I10").Value = '_PMC_elsevier'

$ws.Range("E11").Value = '[Wanbo%Zhu%NULL%0,    Kai%Xie%NULL%0,    Kai%Xie%NULL%0,    Hui%Lu%NULL%0,    Lei%Xu%bayinhexl@126.com%0,    Shusheng%Zhou%zhouss108@163.com%0,    Shiyuan%Fang%fangshiyuan2008@126.com%0]'
$ws.Range("F11").Value = 'PMC7194021'
$ws.Range("I11").Value = '_PMC'

$ws.Range("E12").Value = '[Anthony%Dangis%NULL%0,    Christopher%Gieraerts%NULL%2,    Christopher%Gieraerts%NULL%0,    Yves%De Bruecker%NULL%2,    Yves%De Bruecker%NULL%0,    Lode%Janssen%NULL%2,    Lode%Janssen%NULL%0,    Hanne%Valgaeren%NULL%2,    Hanne%Valgaeren%NULL%0,    Dagmar%Obbels%NULL%2,    Dagmar%Obbels%NULL%0,    Marc%Gillis%NULL%2,    Marc%Gillis%NULL%0,    Marc%Van Ranst%NULL%0,    Marc%Van Ranst%NULL%0,    Johan%Frans%NULL%1,    Annick%Demeyere%NULL%2,    Annick%Demeyere%NULL%0,    Rolf%Symons%rolf.symons@imelda.be%2,    Rolf%Symons%rolf.symons@imelda.be%0]'
$ws.Range("I12").Value = '_PMC'

$ws.Range("C13").Value = 'Accuracy and Reproducibility of Low-Dose Submillisievert Chest CT for the Diagnosis of COVID-19'
$ws.Range("D13").Value = 'Purpose
To demonstrate the accuracy and reproducibility of low-dose submillisievert chest CT for the diagnosis of coronavirus disease 2019 (COVID-19) infection in patients in the emergency department.


Materials and Methods
This was a Health Insurance Portability and Accountability Act–compliant, institutional review board–approved retrospective study.

 From March 14 to 24, 2020, 192 patients in the emergency department with symptoms suggestive of COVID-19 infection were studied by using low-dose chest CT and real-time reverse transcription polymerase chain reaction (RT-PCR).

 Image analysis included the likelihood of COVID-19 infection and the semiquantitative extent of lung involvement.

 CT images were analyzed by two radiologists blinded to the RT-PCR results.

 Reproducibility was assessed using the McNemar test and intraclass correlation coefficient.

 Time between CT acquisition and report was measured.


Results
When compared with RT-PCR, low-dose submillisievert chest CT demonstrated excellent sensitivity, specificity, positive predictive value, negative predictive value, and accuracy for diagnosis of COVID-19 (86.7%, 93.6%, 91.1%, 90.3%, and 90.2%, respectively), in particular in patients with clinical symptoms for more than 48 hours (95.6%, 93.2%, 91.5%, 96.5%, and 94.4%, respectively).

 In patients with a positive CT result, the likelihood of disease increased from 43.2% (pretest probability) to 91.1% or 91.4% (posttest probability), while in patients with a negative CT result, the likelihood of disease decreased to 9.6% or 3.7% for all patients or those with clinical symptoms for &gt;48 hours.

 The prevalence of alternative diagnoses based on chest CT in patients without COVID-19 infection was 17.6%.

 The mean effective radiation dose was 0.56 mSv ± 0.25 (standard deviation).

 Median time between CT acquisition and report was 25 minutes (interquartile range: 13–49 minutes).

 Intra- and interreader reproducibility of CT was excellent (all intraclass correlation coefficients ≥ 0.95) without significant bias in the Bland-Altman analysis.


Conclusion
Low-dose submillisievert chest CT allows for rapid, accurate, and reproducible assessment of COVID-19 infection in patients in the emergency department, in particular in patients with symptoms lasting longer than 48 hours.

 Chest CT has the additional advantage of offering alternative diagnoses in a significant subset of patients.


'
$ws.Range("E13").Value = '[Anthony%Dangis%NULL%0, Christopher%Gieraerts%NULL%2, Christopher%Gieraerts%NULL%0, Yves%De Bruecker%NULL%2, Yves%De Bruecker%NULL%0, Lode%Janssen%NULL%2, Lode%Janssen%NULL%0, Hanne%Valgaeren%NULL%2, Hanne%Valgaeren%NULL%0, Dagmar%Obbels%NULL%2, Dagmar%Obbels%NULL%0, Marc%Gillis%NULL%2, Marc%Gillis%NULL%0, Marc%Van Ranst%NULL%2, Marc%Van Ranst%NULL%0, Johan%Frans%NULL%1, Annick%Demeyere%NULL%2, Annick%Demeyere%NULL%0, Rolf%Symons%rolf.symons@imelda.be%2, Rolf%Symons%rolf.symons@imelda.be%0]'
$ws.Range("F13").Value = 'PMC7233439'
$ws.Range("G13").Value = 'PMC'
$ws.Range("H13").Value = '2020-04-16'
$ws.Range("I13").Value = ''
